{"js": "const NEW_VALUES = [\"4+84=\", \"83-43=\", \"45-23=\", \"89-81=\", \"10+31=\", \"15+38=\", \"69-38=\", \"88-19=\", \"1+26=\", \"94-20=\", \"15+55=\", \"39+33=\", \"71+8=\", \"65+30=\", \"75-69=\", \"66-53=\", \"54-38=\", \"49-23=\", \"94-77=\", \"41-34=\", \"89-24=\", \"19+14=\", \"52-9=\", \"35+56=\", \"34-6=\", \"27-3=\", \"95+1=\", \"88-5=\", \"48-29=\", \"12+73=\", \"35+16=\", \"38-30=\", \"96-16=\", \"57-17=\", \"24-7=\", \"86-48=\", \"49+43=\", \"32+35=\", \"79-54=\", \"62+30=\", \"99-25=\", \"13+48=\", \"46-2=\", \"86-44=\", \"56-38=\", \"2+84=\", \"85-8=\", \"37+52=\", \"3+10=\", \"19-10=\", \"34+32=\", \"72+2=\", \"44-1=\", \"28-9=\", \"65-2=\", \"15+63=\", \"59-33=\", \"50-17=\", \"38+7=\", \"30-5=\", \"58-17=\", \"32-2=\", \"85-84=\", \"18+40=\", \"28-22=\", \"35+20=\", \"23+68=\", \"23+61=\", \"8+8=\", \"86-63=\", \"88-80=\", \"62+10=\", \"4+30=\", \"50+46=\", \"72-35=\", \"15+74=\", \"35+52=\", \"15+7=\", \"7+26=\", \"54+18=\", \"52+34=\", \"53-30=\", \"4+78=\", \"52-45=\", \"53-22=\", \"22+69=\", \"2+90=\", \"63-14=\", \"91-6=\", \"90-40=\", \"48-35=\", \"85-33=\", \"49+27=\", \"90-81=\", \"92-62=\", \"17-1=\", \"24+18=\", \"60+33=\", \"37+48=\", \"18+43=\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst COLS = 5;\nconst rowCount = table.rowCount;\n\nlet k = 0;\nfor (let r = 0; r < rowCount && k < NEW_VALUES.length; r++) {\n  for (let c = 0; c < COLS && k < NEW_VALUES.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = NEW_VALUES[k];\n    k++;\n  }\n}\n\nawait context.sync();\n", "ps1": "$NewValues = @(\"4+84=\",\"83-43=\",\"45-23=\",\"89-81=\",\"10+31=\",\"15+38=\",\"69-38=\",\"88-19=\",\"1+26=\",\"94-20=\",\"15+55=\",\"39+33=\",\"71+8=\",\"65+30=\",\"75-69=\",\"66-53=\",\"54-38=\",\"49-23=\",\"94-77=\",\"41-34=\",\"89-24=\",\"19+14=\",\"52-9=\",\"35+56=\",\"34-6=\",\"27-3=\",\"95+1=\",\"88-5=\",\"48-29=\",\"12+73=\",\"35+16=\",\"38-30=\",\"96-16=\",\"57-17=\",\"24-7=\",\"86-48=\",\"49+43=\",\"32+35=\",\"79-54=\",\"62+30=\",\"99-25=\",\"13+48=\",\"46-2=\",\"86-44=\",\"56-38=\",\"2+84=\",\"85-8=\",\"37+52=\",\"3+10=\",\"19-10=\",\"34+32=\",\"72+2=\",\"44-1=\",\"28-9=\",\"65-2=\",\"15+63=\",\"59-33=\",\"50-17=\",\"38+7=\",\"30-5=\",\"58-17=\",\"32-2=\",\"85-84=\",\"18+40=\",\"28-22=\",\"35+20=\",\"23+68=\",\"23+61=\",\"8+8=\",\"86-63=\",\"88-80=\",\"62+10=\",\"4+30=\",\"50+46=\",\"72-35=\",\"15+74=\",\"35+52=\",\"15+7=\",\"7+26=\",\"54+18=\",\"52+34=\",\"53-30=\",\"4+78=\",\"52-45=\",\"53-22=\",\"22+69=\",\"2+90=\",\"63-14=\",\"91-6=\",\"90-40=\",\"48-35=\",\"85-33=\",\"49+27=\",\"90-81=\",\"92-62=\",\"17-1=\",\"24+18=\",\"60+33=\",\"37+48=\",\"18+43=\")\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$cols = 5\n$rows = $tbl.Rows.Count\n\n$k = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        if ($k -ge $NewValues.Count) { break }\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $NewValues[$k]\n        $k++\n    }\n}\n"}
